# Update the public EPEX Spot prices workbook.
$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert a new date column (08-dec) before the
#     01-oct. block, shifting the 01-oct. .. 31-oct. columns one column
#     to the right (EE -> FJ becomes the new extent).
$wsSpot = $wb.Worksheets.Item("Prix Spot")
$wsSpot.Columns("EE:EE").Insert()
$wsSpot.Range("EE1").Value = "08-dec"
$wsSpot.Range("EE2:EE25").Value = "-"

# --- Sheet "Gaz": append the two new daily quotes.
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A163:A164").NumberFormat = "@"
$wsGaz.Range("A163").Value = "2025-12-06"
$wsGaz.Range("A164").Value = "2025-12-07"
$wsGaz.Range("A163:A164").Style = "Normal"
$wsGaz.Range("B163").Value = 25.905
$wsGaz.Range("B164").Value = 25.905

# --- Sheet "CO2": append the two new daily quotes.
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A164:A165").NumberFormat = "@"
$wsCo2.Range("A164").Value = "2025-12-06"
$wsCo2.Range("A165").Value = "2025-12-07"
$wsCo2.Range("A164:A165").Style = "Normal"
$wsCo2.Range("B164").Value = 81.78
$wsCo2.Range("B165").Value = 81.78
